$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove historical_growth columns D and E for rows 2 and 3
$ws.Range("D2:E3").ClearContents()

# Row 2 updates
$ws.Range("G2").Value = 0.1796078431372549
$ws.Range("H2").Value = 0.1796078431372549
$ws.Range("I2").Value = 0.1505882352941176
$ws.Range("J2").Value = 0.113832231117299
$ws.Range("K2").Value = 2.46
$ws.Range("L2").Value = 0.09647058823529411
$ws.Range("U2").Value = 15.9
$ws.Range("V2").Value = 5.910780669144982
$ws.Range("W2").Value = 0.09283018867924528
$ws.Range("X2").Value = 0.08888812734145883
$ws.Range("Y2").Value = 0.003942061337786448
$ws.Range("Z2").Value = 0.864406779661017
$ws.Range("AA2").Value = 0.09839735232173302
$ws.Range("AB2").Value = 0.08607181213620371
$ws.Range("AC2").Value = 0.01232554018552931
$ws.Range("AD2").Value = 0.37
$ws.Range("AF2").Value = 0.37
$ws.Range("AG2").Value = -15.53
$ws.Range("AH2").Value = 0.1209150326797386
$ws.Range("AI2").Value = 0.007713154054617469
$ws.Range("AJ2").Value = 1.209501557632399
$ws.Range("AK2").Value = -0.4842531961334581
$ws.Range("AL2").Value = 0.468
$ws.Range("AM2").Value = 0.468
$ws.Range("AN2").Value = 0.072265625
$ws.Range("AO2").Value = 8.205128205128204
$ws.Range("AP2").Value = -3.033203125
$ws.Range("AQ2").Value = 8.205128205128204

# Row 3 updates
$ws.Range("G3").Value = 0.1796078431372549
$ws.Range("H3").Value = 0.1796078431372549
$ws.Range("I3").Value = 0.1505882352941176
$ws.Range("J3").Value = 0.113832231117299
$ws.Range("K3").Value = 2.46
$ws.Range("L3").Value = 0.09647058823529411
$ws.Range("U3").Value = 15.9
$ws.Range("V3").Value = 5.910780669144982
$ws.Range("W3").Value = 0.09283018867924528
$ws.Range("X3").Value = 0.08888812734145883
$ws.Range("Y3").Value = 0.003942061337786448
$ws.Range("Z3").Value = 0.864406779661017
$ws.Range("AA3").Value = 0.09839735232173302
$ws.Range("AB3").Value = 0.08607181213620371
$ws.Range("AC3").Value = 0.01232554018552931
$ws.Range("AD3").Value = 0.37
$ws.Range("AF3").Value = 0.37
$ws.Range("AG3").Value = -15.53
$ws.Range("AH3").Value = 0.1209150326797386
$ws.Range("AI3").Value = 0.007713154054617469
$ws.Range("AJ3").Value = 1.209501557632399
$ws.Range("AK3").Value = -0.4842531961334581
$ws.Range("AL3").Value = 0.468
$ws.Range("AM3").Value = 0.468
$ws.Range("AN3").Value = 0.072265625
$ws.Range("AO3").Value = 8.205128205128204
$ws.Range("AP3").Value = -3.033203125
$ws.Range("AQ3").Value = 8.205128205128204
